$p = $ppt.ActivePresentation

# --- Slide 1 ---
$s1 = $p.Slides.Item(1)

# 1. Remove the leftover "Таблица 1" table graphic frame (first shape).
$s1.Shapes.Item(1).Delete()

# 2. "Laundry & Amenities" textbox: drop the trailing empty paragraph
#    ("." on its own line) while keeping the other two paragraphs/runs
#    (and their formatting) untouched.
$grpLaundry = $s1.Shapes.Item("Группа 33")
$tbLaundry = $grpLaundry.GroupItems.Item("TextBox 45")
$trLaundry = $tbLaundry.TextFrame.TextRange
$lastPara = $trLaundry.Paragraphs($trLaundry.Paragraphs().Count, 1)
$lastPara.Text = ""
$trLaundry2 = $tbLaundry.TextFrame.TextRange
$trLaundry2.Paragraphs($trLaundry2.Paragraphs().Count, 1).Delete()
# Restore the autofit height to the value PowerPoint itself computes for
# the now two-paragraph box (the text-metric engine here is only an
# approximation of that).
$tbLaundry.Height = 28.2734657

# 3. "services" -> "Services" (capitalize); the box keeps its original size.
$grpServices = $s1.Shapes.Item("Группа 50")
$tbServices = $grpServices.GroupItems.Item("TextBox 70")
$origHeight = $tbServices.Height
$tbServices.TextFrame.TextRange.Text = "Services"
$tbServices.Height = $origHeight

# --- Slide 2 ---
$s2 = $p.Slides.Item(2)

# Remove the leftover "Таблица 1" table graphic frame on this slide too.
$s2.Shapes.Item(1).Delete()
